$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 11.921572
$ws.Cells.Item(2, 8).Value = 35.764716
$ws.Cells.Item(2, 9).Value = 0.006971694289596158
$ws.Cells.Item(2, 10).Value = 0.006971694289596159
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 21.08181366666667
$ws.Cells.Item(2, 14).Value = 63.245441
$ws.Cells.Item(2, 15).Value = 0.0571606014598545
$ws.Cells.Item(2, 16).Value = 0.0571606014598545
$ws.Cells.Item(2, 17).Value = 251.3283595177506
$ws.Cells.Item(2, 18).Value = 2261.955235659756
$ws.Cells.Item(2, 19).Value = 0.0003985062387875494
$ws.Cells.Item(2, 20).Value = 0.0003985062387875495

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 11.921572
$ws.Cells.Item(3, 8).Value = 35.764716
$ws.Cells.Item(3, 9).Value = 0.006971694289596158
$ws.Cells.Item(3, 10).Value = 0.006971694289596159
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 301.6001486666667
$ws.Cells.Item(3, 14).Value = 904.800446
$ws.Cells.Item(3, 15).Value = 0.8177496571571792
$ws.Cells.Item(3, 16).Value = 0.8177496571571792
$ws.Cells.Item(3, 17).Value = 3595.547887540371
$ws.Cells.Item(3, 18).Value = 32359.93098786334
$ws.Cells.Item(3, 19).Value = 0.005701100615121922
$ws.Cells.Item(3, 20).Value = 0.005701100615121923

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 11.921572
$ws.Cells.Item(4, 8).Value = 35.764716
$ws.Cells.Item(4, 9).Value = 0.006971694289596158
$ws.Cells.Item(4, 10).Value = 0.006971694289596159
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 46.13524966666667
$ws.Cells.Item(4, 14).Value = 138.405749
$ws.Cells.Item(4, 15).Value = 0.1250897413829664
$ws.Cells.Item(4, 16).Value = 0.1250897413829664
$ws.Cells.Item(4, 17).Value = 550.0047006391427
$ws.Cells.Item(4, 18).Value = 4950.042305752285
$ws.Cells.Item(4, 19).Value = 0.0008720874356866868
$ws.Cells.Item(4, 20).Value = 0.0008720874356866869

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1516.953124666667
$ws.Cells.Item(5, 8).Value = 4550.859374
$ws.Cells.Item(5, 9).Value = 0.8871089682487887
$ws.Cells.Item(5, 10).Value = 0.8871089682487888
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 21.08181366666667
$ws.Cells.Item(5, 14).Value = 63.245441
$ws.Cells.Item(5, 15).Value = 0.0571606014598545
$ws.Cells.Item(5, 16).Value = 0.0571606014598545
$ws.Cells.Item(5, 17).Value = 31980.12311529043
$ws.Cells.Item(5, 18).Value = 287821.1080376139
$ws.Cells.Item(5, 19).Value = 0.05070768218553173
$ws.Cells.Item(5, 20).Value = 0.05070768218553174

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1516.953124666667
$ws.Cells.Item(6, 8).Value = 4550.859374
$ws.Cells.Item(6, 9).Value = 0.8871089682487887
$ws.Cells.Item(6, 10).Value = 0.8871089682487888
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 301.6001486666667
$ws.Cells.Item(6, 14).Value = 904.800446
$ws.Cells.Item(6, 15).Value = 0.8177496571571792
$ws.Cells.Item(6, 16).Value = 0.8177496571571792
$ws.Cells.Item(6, 17).Value = 457513.2879198312
$ws.Cells.Item(6, 18).Value = 4117619.59127848
$ws.Cells.Item(6, 19).Value = 0.7254330546465059
$ws.Cells.Item(6, 20).Value = 0.725433054646506

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1516.953124666667
$ws.Cells.Item(7, 8).Value = 4550.859374
$ws.Cells.Item(7, 9).Value = 0.8871089682487887
$ws.Cells.Item(7, 10).Value = 0.8871089682487888
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 46.13524966666667
$ws.Cells.Item(7, 14).Value = 138.405749
$ws.Cells.Item(7, 15).Value = 0.1250897413829664
$ws.Cells.Item(7, 16).Value = 0.1250897413829664
$ws.Cells.Item(7, 17).Value = 69985.0111391268
$ws.Cells.Item(7, 18).Value = 629865.1002521411
$ws.Cells.Item(7, 19).Value = 0.1109682314167511
$ws.Cells.Item(7, 20).Value = 0.1109682314167511

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 181.1216836666667
$ws.Cells.Item(8, 8).Value = 543.365051
$ws.Cells.Item(8, 9).Value = 0.1059193374616151
$ws.Cells.Item(8, 10).Value = 0.1059193374616151
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 21.08181366666667
$ws.Cells.Item(8, 14).Value = 63.245441
$ws.Cells.Item(8, 15).Value = 0.0571606014598545
$ws.Cells.Item(8, 16).Value = 0.0571606014598545
$ws.Cells.Item(8, 17).Value = 3818.373586053609
$ws.Cells.Item(8, 18).Value = 34365.36227448249
$ws.Cells.Item(8, 19).Value = 0.006054413035535218
$ws.Cells.Item(8, 20).Value = 0.006054413035535219

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 181.1216836666667
$ws.Cells.Item(9, 8).Value = 543.365051
$ws.Cells.Item(9, 9).Value = 0.1059193374616151
$ws.Cells.Item(9, 10).Value = 0.1059193374616151
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 301.6001486666667
$ws.Cells.Item(9, 14).Value = 904.800446
$ws.Cells.Item(9, 15).Value = 0.8177496571571792
$ws.Cells.Item(9, 16).Value = 0.8177496571571792
$ws.Cells.Item(9, 17).Value = 54626.32672062363
$ws.Cells.Item(9, 18).Value = 491636.9404856127
$ws.Cells.Item(9, 19).Value = 0.08661550189555133
$ws.Cells.Item(9, 20).Value = 0.08661550189555134

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 181.1216836666667
$ws.Cells.Item(10, 8).Value = 543.365051
$ws.Cells.Item(10, 9).Value = 0.1059193374616151
$ws.Cells.Item(10, 10).Value = 0.1059193374616151
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 46.13524966666667
$ws.Cells.Item(10, 14).Value = 138.405749
$ws.Cells.Item(10, 15).Value = 0.1250897413829664
$ws.Cells.Item(10, 16).Value = 0.1250897413829664
$ws.Cells.Item(10, 17).Value = 8356.09409600869
$ws.Cells.Item(10, 18).Value = 75204.84686407821
$ws.Cells.Item(10, 19).Value = 0.01324942253052857
$ws.Cells.Item(10, 20).Value = 0.01324942253052858

